$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test_1")

# Remove the Gas_flow, Water_flow and Oil_flow columns (D:F) entirely,
# shifting Pressure_atm / Volumetric flow_m3/h / Water / Hydrogen sulfide left.
$ws.Range("D1:F1").EntireColumn.Delete()

# Update the active selection to reflect where the user ended up after the edit.
$ws.Range("C5").Select()
